# Add bold markers ("**...**") around the section-header labels inside the
# "Razão da Falha" guidance text in C2 (Sheet1), matching the source edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
**Razão da Falha**		
	Uso incorreto do modelo de cumprimento (ex: BraHello	
	 quando há histórico de Chatbot).	

**Comportamento Correto**:		
	BraHello>	
	Sem interação prévia com Chatbot.	
	BraAcceptTransfer	>
	Quando há conversa prévia com Chatbot.	

**Comportamento Incorreto**:		
	Não usar modelo específico ou usar modelo inadequado ao contexto.	

**Exemplos**:		
	**Correto:	
		Cenário: Cliente relata depósito não creditado após interagir com Chatbot.
		Ação do Agente: Usa 
		IngAcceptTransfer
		 e menciona: "Vi que seu contato é sobre o depósito".
	**Incorreto**:	
		Cenário: Histórico de Chatbot visível.
		Ação do Agente: Usa 
		BraHello
		 ou mensagem genérica ("Como posso ajudar?").

**Notas**:		
	Personalização (ex: incluir nome do cliente) é permitida, mas não substitui o modelo obrigatório.	
	E-mails: Selecionar modelo de e-mail correspondente à consulta.	
	Atraso > 60 segundos na abertura é avaliado em Questão 9 (Tempos de Espera).	
"@

$ws.Range("C2").Value = $newText
